$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while forcing it to stay a text
# string even when the text looks like a number (e.g. "591.90"), and then
# strip the temporary "Text" number format again so the cell ends up with
# its original (default) style.
function Set-TextCell($worksheet, $row, $col, $text) {
    $cell = $worksheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$sub3 = [string][char]0x2083

# Row -> new Price (column D) / new Volume(1h) (column E) values.
# A $null Price means that column is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "61.147.15";            E = "  +0.62%  " },
    @{ Row = 3;  D = "2.930.06";             E = "  +0.78%  " },
    @{ Row = 4;  D = $null;                  E = "  +0.03%  " },
    @{ Row = 5;  D = "591.90";               E = "  +0.83%  " },
    @{ Row = 6;  D = "145.47";               E = "  +0.89%  " },
    @{ Row = 8;  D = "0.505";                E = "  +0.62%  " },
    @{ Row = 9;  D = "6.98";                 E = "  +4.07%  " },
    @{ Row = 10; D = "0.143";                E = "  -0.02%  " },
    @{ Row = 11; D = "0.441";                E = "  -0.54%  " },
    @{ Row = 12; D = "0.0000225";            E = "  +0.18%  " },
    @{ Row = 13; D = "33.77";                E = "  +1.11%  " },
    @{ Row = 14; D = $null;                  E = "  -0.22%  " },
    @{ Row = 15; D = "3.419.26";             E = "  +0.94%  " },
    @{ Row = 16; D = "60.989.72";            E = "  +0.41%  " },
    @{ Row = 17; D = "6.73";                 E = "  +0.89%  " },
    @{ Row = 18; D = "2.932.10";             E = "  +0.86%  " },
    @{ Row = 19; D = "437.14";               E = "  +2.21%  " },
    @{ Row = 20; D = "13.44";                E = "  -0.51%  " },
    @{ Row = 21; D = "0.679";                E = "  -0.34%  " },
    @{ Row = 22; D = "7.11";                 E = "  +0.78%  " },
    @{ Row = 23; D = "81.63";                E = "  +1.11%  " },
    @{ Row = 24; D = "11.01";                E = "  +1.76%  " },
    @{ Row = 25; D = $null;                  E = "  -0.23%  " },
    @{ Row = 26; D = "11.90";                E = "  +0.33%  " },
    @{ Row = 27; D = $null;                  E = "  +0.05%  " },
    @{ Row = 28; D = "2.26";                 E = "  +3.17%  " },
    @{ Row = 29; D = $null;                  E = "  +0.31%  " },
    @{ Row = 30; D = "7.01";                 E = "  -2.41%  " },
    @{ Row = 31; D = "0.110";                E = "  +3.40%  " },
    @{ Row = 32; D = "26.66";                E = "  +1.06%  " },
    @{ Row = 33; D = $null;                  E = "  +0.10%  " },
    @{ Row = 34; D = "0.0" + $sub3 + "0870"; E = "  +1.58%  " },
    @{ Row = 35; D = "1.01";                 E = "  +0.46%  " },
    @{ Row = 36; D = "5.64";                 E = "  +1.28%  " },
    @{ Row = 37; D = "3.01";                 E = "  +0.00%  " },
    @{ Row = 38; D = $null;                  E = "  +0.37%  " },
    @{ Row = 39; D = "1.99";                 E = "  +0.60%  " },
    @{ Row = 40; D = "8.60";                 E = "  +0.39%  " },
    @{ Row = 41; D = "42.12";                E = "  +2.02%  " },
    @{ Row = 42; D = $null;                  E = "  -2.38%  " },
    @{ Row = 43; D = "376.35";               E = "  +0.66%  " },
    @{ Row = 44; D = "0.0347";               E = "  -0.91%  " },
    @{ Row = 45; D = "2.694.58";             E = "  +0.05%  " },
    @{ Row = 46; D = "133.05";               E = "  +0.54%  " },
    @{ Row = 48; D = "23.99";                E = "  -0.24%  " },
    @{ Row = 49; D = $null;                  E = "  -0.21%  " },
    @{ Row = 50; D = $null;                  E = "  -1.36%  " },
    @{ Row = 51; D = "0.125";                E = "  +0.99%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextCell $ws $r 4 $u.D
    }
    Set-TextCell $ws $r 5 $u.E
}
